$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit reorders the data rows (2-25) of the sheet: every row keeps its
# full set of column values (A:T) together, but rows are re-sequenced into
# a new order. Map: new row number -> old row number it takes its data from.
$rowMap = @{ 2 = 6; 3 = 13; 4 = 20; 5 = 4; 6 = 14; 7 = 24; 8 = 16; 9 = 10; 10 = 15; 11 = 9; 12 = 19; 13 = 5; 14 = 2; 15 = 3; 16 = 18; 17 = 7; 18 = 25; 19 = 22; 20 = 23; 21 = 11; 22 = 17; 23 = 12; 24 = 21; 25 = 8 }

# Snapshot every source row's A:T values before any writes happen, so that
# writing a new row never clobbers data that is still needed as a source.
$snapshot = @{}
foreach ($newRow in $rowMap.Keys) {
    $oldRow = $rowMap[$newRow]
    if (-not $snapshot.ContainsKey($oldRow)) {
        $snapshot[$oldRow] = $ws.Range("A" + $oldRow + ":T" + $oldRow).Value()
    }
}

foreach ($newRow in ($rowMap.Keys | Sort-Object)) {
    $oldRow = $rowMap[$newRow]
    $ws.Range("A" + $newRow + ":T" + $newRow).Value = $snapshot[$oldRow]
}
